# Updated to include parts for Blinker Circuit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section marker row ---
$ws.Range("A5").Value = "x"

# --- New parts data (rows 6-13) ---
$partRows = @(
    @{ Row=6;  Item=4;  Desc="8-DIP Socket for 555";  Part="A08-LC-TT";         Qty=10; Price=0.26;  Url="http://ca.mouser.com/ProductDetail/TE-Connectivity/A08-LC-TT" },
    @{ Row=7;  Item=5;  Desc="555 Timer Circuit";      Part="LM555CN";           Qty=10; Price=0.59;  Url="http://ca.mouser.com/ProductDetail/Texas-Instruments/LM555CN" },
    @{ Row=8;  Item=6;  Desc="5V Voltage Regulator";   Part="L7805ABD2T-TR";     Qty=10; Price=0.92;  Url="http://ca.mouser.com/ProductDetail/STMicroelectronics/L7805ABD2T-TR" },
    @{ Row=9;  Item=7;  Desc="10uF Capacitor";         Part="TMK316F106ZL-T";    Qty=10; Price=0.3;   Url="http://ca.mouser.com/ProductDetail/Taiyo-Yuden/TMK316F106ZL-T" },
    @{ Row=10; Item=8;  Desc="10nF Capacitor";         Part="CC1206KRX7R9BB103"; Qty=10; Price=0.19;  Url="http://ca.mouser.com/ProductDetail/Yageo/CC1206KRX7R9BB103" },
    @{ Row=11; Item=9;  Desc="1k Resistor";             Part="ERJ-14YJ102U";      Qty=10; Price=0.32;  Url="http://ca.mouser.com/ProductDetail/Panasonic/ERJ-14YJ102U" },
    @{ Row=12; Item=10; Desc="22k Resistor";            Part="ERJ-14YJ223U";      Qty=10; Price=0.26;  Url="http://ca.mouser.com/ProductDetail/Panasonic/ERJ-14YJ223U" },
    @{ Row=13; Item=11; Desc="50k Resistor";            Part="CRCW121049K9FKEA";  Qty=10; Price=0.441; Url="http://ca.mouser.com/ProductDetail/Vishay/CRCW121049K9FKEA" }
)

foreach ($r in $partRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Item
    $ws.Cells.Item($row, 2).Value = $r.Desc
    $ws.Cells.Item($row, 3).Value = $r.Part
    $ws.Cells.Item($row, 4).Value = $r.Qty
    $ws.Cells.Item($row, 5).Value = $r.Price
    $ws.Cells.Item($row, 6).Value = "CAD"
    $ws.Cells.Item($row, 7).Formula = "=D" + $row + "*E" + $row

    # Hyperlink on the part-number cell, matching existing Hyperlink style.
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 3), $r.Url)
    $ws.Cells.Item($row, 3).Style = $ws.Range("C2").Style
}

# --- Grand-total formula next to the header row ---
$ws.Range("H1").Formula = "=SUM(G2:G30)"

# --- Resize the table / autofilter to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G13"))

# --- View state tweaks ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("H2").Select() | Out-Null
